# "easy view finish" — add the missing "Close" control (rounded-rect button
# + "Close" textbox, grouped) to slide 1 of the management-dashboard
# interface prototype. The new group is a sibling of the existing
# "组合 23" (group id 24) close-control group, just nudged to a slightly
# different position, so we clone that existing group (to inherit its
# exact shape style / text formatting) and then rename + reposition the
# clone to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the existing "close control" group (p:cNvPr id="18" ... up the
# chain; the one we want to clone is id=24, name "组合 23").
$source = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 24) {
        $source = $candidate
    }
}

# Clone it — this preserves the full p:style / lstStyle / effect
# formatting of both child shapes (the rounded-rect button and the
# "Close" textbox) exactly, which plain AddShape/AddTextbox calls cannot
# reproduce.
$dup = $source.Duplicate()
$newGroup = $dup.Item(1)

# Rename the new group and its "Close" textbox to fresh default-style
# names (the duplicated rounded rectangle naturally keeps its source
# name, "矩形: 圆角 24", matching the target).
$newGroup.Name = "组合 37"
$groupItems = $newGroup.GroupItems
$groupItems.Item(2).Name = "文本框 39"

# Move the clone to its target position (EMU -> points, 1 pt = 12700 EMU;
# rounded to 4 decimal places so the float32 COM round-trip lands on the
# exact target EMU values: off x="10581730" y="5025245").
$newGroup.Left = 833.2071
$newGroup.Top = 395.6886
